$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title text update (October -> November) ---
$ws.Range("A2").Value = "by State, by Sector, Year-to-Date through November 2016 and 2015 (Thousand Megawatthours)"

# --- Column header row 6: October -> November for all YTD header cells ---
$ws.Range("B6").Value = "November 2016 YTD"
$ws.Range("C6").Value = "November 2015 YTD"
$ws.Range("E6").Value = "November 2016 YTD"
$ws.Range("F6").Value = "November 2015 YTD"
$ws.Range("G6").Value = "November 2016 YTD"
$ws.Range("H6").Value = "November 2015 YTD"
$ws.Range("I6").Value = "November 2016 YTD"
$ws.Range("J6").Value = "November 2015 YTD"
$ws.Range("K6").Value = "November 2016 YTD"
$ws.Range("L6").Value = "November 2015 YTD"

# --- Data cell updates ---
$ws.Range("B7").Value = 29892
$ws.Range("C7").Value = 28903
$ws.Range("D7").Value = 0.034
$ws.Range("G7").Value = 29892
$ws.Range("H7").Value = 28903

$ws.Range("B8").Value = 15022
$ws.Range("C8").Value = 15854
$ws.Range("D8").Value = -0.053
$ws.Range("G8").Value = 15022
$ws.Range("H8").Value = 15854

$ws.Range("B10").Value = 5039
$ws.Range("C10").Value = 4493
$ws.Range("D10").Value = 0.121
$ws.Range("G10").Value = 5039
$ws.Range("H10").Value = 4493

$ws.Range("B11").Value = 9831
$ws.Range("C11").Value = 8555
$ws.Range("D11").Value = 0.149
$ws.Range("G11").Value = 9831
$ws.Range("H11").Value = 8555

$ws.Range("B14").Value = 140359
$ws.Range("C14").Value = 144192
$ws.Range("D14").Value = -0.027
$ws.Range("G14").Value = 140359
$ws.Range("H14").Value = 144192

$ws.Range("B15").Value = 26839
$ws.Range("C15").Value = 30211
$ws.Range("D15").Value = -0.112
$ws.Range("G15").Value = 26839
$ws.Range("H15").Value = 30211

$ws.Range("B16").Value = 38017
$ws.Range("C16").Value = 40729
$ws.Range("D16").Value = -0.067
$ws.Range("G16").Value = 38017
$ws.Range("H16").Value = 40729

$ws.Range("B17").Value = 75504
$ws.Range("C17").Value = 73251
$ws.Range("D17").Value = 0.031
$ws.Range("G17").Value = 75504
$ws.Range("H17").Value = 73251

$ws.Range("B18").Value = 143113
$ws.Range("C18").Value = 140516
$ws.Range("D18").Value = 0.018
$ws.Range("E18").Value = 22845
$ws.Range("F18").Value = 21347
$ws.Range("G18").Value = 120268
$ws.Range("H18").Value = 119169

$ws.Range("B19").Value = 89380
$ws.Range("C19").Value = 88569
$ws.Range("D19").Value = 0.009
$ws.Range("G19").Value = 89380
$ws.Range("H19").Value = 88569

$ws.Range("B21").Value = 29285
$ws.Range("C21").Value = 27061
$ws.Range("D21").Value = 0.082
$ws.Range("E21").Value = 22845
$ws.Range("F21").Value = 21347
$ws.Range("G21").Value = 6440
$ws.Range("H21").Value = 5714

$ws.Range("B22").Value = 15193
$ws.Range("C22").Value = 15747
$ws.Range("D22").Value = -0.035
$ws.Range("G22").Value = 15193
$ws.Range("H22").Value = 15747

$ws.Range("B23").Value = 9255
$ws.Range("C23").Value = 9137
$ws.Range("D23").Value = 0.013
$ws.Range("G23").Value = 9255
$ws.Range("H23").Value = 9137

$ws.Range("B24").Value = 41372
$ws.Range("C24").Value = 42460
$ws.Range("D24").Value = -0.026
$ws.Range("E24").Value = 37128
$ws.Range("F24").Value = 37673
$ws.Range("G24").Value = 4243
$ws.Range("H24").Value = 4787

$ws.Range("B25").Value = 4243
$ws.Range("C25").Value = 4787
$ws.Range("D25").Value = -0.114
$ws.Range("G25").Value = 4243
$ws.Range("H25").Value = 4787

$ws.Range("B26").Value = 7333
$ws.Range("C26").Value = 7720
$ws.Range("D26").Value = -0.05
$ws.Range("E26").Value = 7333
$ws.Range("F26").Value = 7720

$ws.Range("B27").Value = 12539
$ws.Range("C27").Value = 11072
$ws.Range("D27").Value = 0.133
$ws.Range("E27").Value = 12539
$ws.Range("F27").Value = 11072

$ws.Range("B28").Value = 8506
$ws.Range("C28").Value = 9521
$ws.Range("D28").Value = -0.107
$ws.Range("E28").Value = 8506
$ws.Range("F28").Value = 9521

$ws.Range("B29").Value = 8749
$ws.Range("C29").Value = 9361
$ws.Range("D29").Value = -0.065
$ws.Range("E29").Value = 8749
$ws.Range("F29").Value = 9361

$ws.Range("B32").Value = 188144
$ws.Range("C32").Value = 182269
$ws.Range("D32").Value = 0.032
$ws.Range("E32").Value = 174679
$ws.Range("F32").Value = 168839
$ws.Range("G32").Value = 13465
$ws.Range("H32").Value = 13431

$ws.Range("B35").Value = 26573
$ws.Range("C35").Value = 25501
$ws.Range("D35").Value = 0.042
$ws.Range("E35").Value = 26573
$ws.Range("F35").Value = 25501

$ws.Range("B36").Value = 31402
$ws.Range("C36").Value = 30777
$ws.Range("D36").Value = 0.02
$ws.Range("E36").Value = 31402
$ws.Range("F36").Value = 30777

$ws.Range("B37").Value = 13465
$ws.Range("C37").Value = 13431
$ws.Range("D37").Value = 0.003
$ws.Range("G37").Value = 13465
$ws.Range("H37").Value = 13431

$ws.Range("B38").Value = 38905
$ws.Range("C38").Value = 38221
$ws.Range("D38").Value = 0.018
$ws.Range("E38").Value = 38905
$ws.Range("F38").Value = 38221

$ws.Range("B39").Value = 50803
$ws.Range("C39").Value = 48760
$ws.Range("D39").Value = 0.042
$ws.Range("E39").Value = 50803
$ws.Range("F39").Value = 48760

$ws.Range("B40").Value = 26996
$ws.Range("C40").Value = 25580
$ws.Range("D40").Value = 0.055
$ws.Range("E40").Value = 26996
$ws.Range("F40").Value = 25580

$ws.Range("B42").Value = 68974
$ws.Range("C42").Value = 72214
$ws.Range("D42").Value = -0.045
$ws.Range("E42").Value = 68974
$ws.Range("F42").Value = 72214

$ws.Range("B43").Value = 36096
$ws.Range("C43").Value = 38381
$ws.Range("D43").Value = -0.06
$ws.Range("E43").Value = 36096
$ws.Range("F43").Value = 38381

$ws.Range("B45").Value = 5911
$ws.Range("C45").Value = 10705
$ws.Range("D45").Value = -0.448
$ws.Range("E45").Value = 5911
$ws.Range("F45").Value = 10705

$ws.Range("B46").Value = 26967
$ws.Range("C46").Value = 23127
$ws.Range("D46").Value = 0.166
$ws.Range("E46").Value = 26967
$ws.Range("F46").Value = 23127

$ws.Range("B47").Value = 66197
$ws.Range("C47").Value = 63227
$ws.Range("D47").Value = 0.047
$ws.Range("E47").Value = 27920
$ws.Range("F47").Value = 26847
$ws.Range("G47").Value = 38277
$ws.Range("H47").Value = 36380

$ws.Range("B48").Value = 12372
$ws.Range("C48").Value = 12530
$ws.Range("D48").Value = -0.013
$ws.Range("E48").Value = 12372
$ws.Range("F48").Value = 12530

$ws.Range("B49").Value = 15548
$ws.Range("C49").Value = 14317
$ws.Range("D49").Value = 0.086
$ws.Range("E49").Value = 15548
$ws.Range("F49").Value = 14317

$ws.Range("B51").Value = 38277
$ws.Range("C51").Value = 36380
$ws.Range("D51").Value = 0.052
$ws.Range("G51").Value = 38277
$ws.Range("H51").Value = 36380

$ws.Range("B52").Value = 29394
$ws.Range("C52").Value = 29540
$ws.Range("D52").Value = -0.005
$ws.Range("E52").Value = 29394
$ws.Range("F52").Value = 29540

$ws.Range("B53").Value = 29394
$ws.Range("C53").Value = 29540
$ws.Range("D53").Value = -0.005
$ws.Range("E53").Value = 29394
$ws.Range("F53").Value = 29540

$ws.Range("B61").Value = 26187
$ws.Range("C61").Value = 24223
$ws.Range("D61").Value = 0.081
$ws.Range("E61").Value = 26187
$ws.Range("F61").Value = 24223

$ws.Range("B62").Value = 17214
$ws.Range("C62").Value = 16922
$ws.Range("D62").Value = 0.017
$ws.Range("E62").Value = 17214
$ws.Range("F62").Value = 16922

$ws.Range("B64").Value = 8973
$ws.Range("C64").Value = 7301
$ws.Range("D64").Value = 0.229
$ws.Range("E64").Value = 8973
$ws.Range("F64").Value = 7301

$ws.Range("B68").Value = 733632
$ws.Range("C68").Value = 727544
$ws.Range("D68").Value = 0.008
$ws.Range("E68").Value = 387127
$ws.Range("F68").Value = 380683
$ws.Range("G68").Value = 346505
$ws.Range("H68").Value = 346861
